# Update cryptocurrency price (D) and volume change (E) cells to reflect latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.728.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.533.96'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.32'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.87%  '
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.751.97'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.532.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('E14').Value = '  -1.54%  '
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.715.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '212.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0681'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.92%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.84%  '
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.52%  '
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('E29').Value = '  -1.23%  '
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('E32').Value = '  +2.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.363.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('E36').Value = '  +2.72%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('E41').Value = '  +6.07%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '62.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.666.26'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₇0974'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0942'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
